# Version3: Minor format modifications
# Target sheet is the second worksheet ("-the-name-of-your- module"),
# which is already the ActiveSheet for this workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content updates -------------------------------------------------

# Row 20 ("Version 3" section header) Run by / Date cell: replace the
# placeholder tester/date text with the real tester name and date.
$ws.Range("G20").Value = "Jarvinia Zhao`n12-03-2023"

# Row 30 "Actual result" cell: "output:" -> "Expected output:"
$ws.Range("D30").Value = "Unable to input needle. haystack should be truncated to BUFFER_SIZE -1 and the rest data will store in needle`nExpected output: Not found"

# Row 32 "Actual result" cell: "Output:" -> "Expected Output:"
$ws.Range("D32").Value = "Can't input the next haystack value`nExpected Output: Not found"

# --- Formatting updates ------------------------------------------------

# B22:B29 ("+ ..." purpose cells in the strstr section) gain a top
# vertical alignment, matching the style already used by the other
# purpose cells (B12:B19, B30:B32) in the sheet.
$ws.Range("B22:B29").VerticalAlignment = -4160

# --- View / selection updates ------------------------------------------

$ws.Activate()
$ws.Range("G33").Select()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
